$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.1122238263487816
$ws.Cells.Item(2, 2).Value = 0.9629124402999878
$ws.Cells.Item(2, 3).Value = 0.01426568906754255
$ws.Cells.Item(2, 4).Value = 0.9982606768608093

$ws.Cells.Item(3, 1).Value = 0.03388161212205887
$ws.Cells.Item(3, 2).Value = 0.9894838333129883
$ws.Cells.Item(3, 3).Value = 0.009175782091915607
$ws.Cells.Item(3, 4).Value = 0.9986085295677185

$ws.Cells.Item(4, 1).Value = 0.02503669448196888
$ws.Cells.Item(4, 2).Value = 0.9899855256080627
$ws.Cells.Item(4, 3).Value = 0.003753812052309513
$ws.Cells.Item(4, 4).Value = 0.9989563822746277

$ws.Cells.Item(5, 1).Value = 0.02060467004776001
$ws.Cells.Item(5, 2).Value = 0.9903668761253357
$ws.Cells.Item(5, 3).Value = 0.002884893910959363
$ws.Cells.Item(5, 4).Value = 0.9989563822746277

$ws.Cells.Item(6, 1).Value = 0.02033344842493534
$ws.Cells.Item(6, 2).Value = 0.9899855256080627
$ws.Cells.Item(6, 3).Value = 0.002187391277402639
$ws.Cells.Item(6, 4).Value = 0.9991883039474487

$ws.Cells.Item(7, 1).Value = 0.01602194644510746
$ws.Cells.Item(7, 2).Value = 0.9932166934013367
$ws.Cells.Item(7, 3).Value = 0.001139729516580701
$ws.Cells.Item(7, 4).Value = 0.9990723729133606

$ws.Cells.Item(8, 1).Value = 0.01216172147542238
$ws.Cells.Item(8, 2).Value = 0.9954643845558167
$ws.Cells.Item(8, 3).Value = 0.0006952740950509906
$ws.Cells.Item(8, 4).Value = 0.9995361566543579

$ws.Cells.Item(9, 1).Value = 0.01009280420839787
$ws.Cells.Item(9, 2).Value = 0.9962671399116516
$ws.Cells.Item(9, 3).Value = 0.001155905425548553
$ws.Cells.Item(9, 4).Value = 0.9995361566543579

$ws.Cells.Item(10, 1).Value = 0.01025610044598579
$ws.Cells.Item(10, 2).Value = 0.9964076280593872
$ws.Cells.Item(10, 3).Value = 0.0005122976726852357
$ws.Cells.Item(10, 4).Value = 0.999768078327179

$ws.Cells.Item(11, 1).Value = 0.01098364219069481
$ws.Cells.Item(11, 2).Value = 0.9956450462341309
$ws.Cells.Item(11, 3).Value = 0.0003453810350038111
$ws.Cells.Item(11, 4).Value = 0.9998840689659119

$ws.Cells.Item(12, 1).Value = 0.010852943174541
$ws.Cells.Item(12, 2).Value = 0.9959259629249573
$ws.Cells.Item(12, 3).Value = 0.0003522520419210196
$ws.Cells.Item(12, 4).Value = 0.999768078327179

$ws.Cells.Item(13, 1).Value = 0.01004441268742085
$ws.Cells.Item(13, 2).Value = 0.9962671399116516
$ws.Cells.Item(13, 3).Value = 0.0005051796906627715
$ws.Cells.Item(13, 4).Value = 0.999768078327179

$ws.Cells.Item(14, 1).Value = 0.01015487499535084
$ws.Cells.Item(14, 2).Value = 0.9959460496902466
$ws.Cells.Item(14, 3).Value = 0.0003993684949818999
$ws.Cells.Item(14, 4).Value = 0.9998840689659119

$ws.Cells.Item(15, 1).Value = 0.0102023659273982
$ws.Cells.Item(15, 2).Value = 0.9958858489990234
$ws.Cells.Item(15, 3).Value = 0.0001680965506238863
$ws.Cells.Item(15, 4).Value = 0.9998840689659119

$ws.Cells.Item(16, 1).Value = 0.01070388313382864
$ws.Cells.Item(16, 2).Value = 0.9957854747772217
$ws.Cells.Item(16, 3).Value = 0.001187594141811132
$ws.Cells.Item(16, 4).Value = 0.9996521472930908

$ws.Cells.Item(17, 1).Value = 0.01119144540280104
$ws.Cells.Item(17, 2).Value = 0.9957252740859985
$ws.Cells.Item(17, 3).Value = 0.0002615421544760466
$ws.Cells.Item(17, 4).Value = 0.9998840689659119

$ws.Cells.Item(18, 1).Value = 0.0108451135456562
$ws.Cells.Item(18, 2).Value = 0.9957252740859985
$ws.Cells.Item(18, 3).Value = 0.0004194485954940319
$ws.Cells.Item(18, 4).Value = 0.999768078327179

$ws.Cells.Item(19, 1).Value = 0.00960271991789341
$ws.Cells.Item(19, 2).Value = 0.9961869120597839
$ws.Cells.Item(19, 3).Value = 0.0002652259427122772
$ws.Cells.Item(19, 4).Value = 0.9998840689659119

$ws.Cells.Item(20, 1).Value = 0.009786729700863361
$ws.Cells.Item(20, 2).Value = 0.9962872266769409
$ws.Cells.Item(20, 3).Value = 0.0001897358597489074
$ws.Cells.Item(20, 4).Value = 0.9998840689659119

$ws.Cells.Item(21, 1).Value = 0.009559962898492813
$ws.Cells.Item(21, 2).Value = 0.9962872266769409
$ws.Cells.Item(21, 3).Value = 0.000239111774135381
$ws.Cells.Item(21, 4).Value = 0.9998840689659119

$ws.Cells.Item(22, 1).Value = 0.009558824822306633
$ws.Cells.Item(22, 2).Value = 0.9963675141334534
$ws.Cells.Item(22, 3).Value = 0.0001472172443754971
$ws.Cells.Item(22, 4).Value = 0.9998840689659119

$ws.Cells.Item(23, 1).Value = 0.009529390372335911
$ws.Cells.Item(23, 2).Value = 0.9961467385292053
$ws.Cells.Item(23, 3).Value = 0.001138524268753827
$ws.Cells.Item(23, 4).Value = 0.9998840689659119

$ws.Cells.Item(24, 1).Value = 0.01111293770372868
$ws.Cells.Item(24, 2).Value = 0.9959661364555359
$ws.Cells.Item(24, 3).Value = 0.00007986863056430593
$ws.Cells.Item(24, 4).Value = 1

$ws.Cells.Item(25, 1).Value = 0.01017675641924143
$ws.Cells.Item(25, 2).Value = 0.9960865378379822
$ws.Cells.Item(25, 3).Value = 0.00006402580766007304
$ws.Cells.Item(25, 4).Value = 1

$ws.Cells.Item(26, 1).Value = 0.009987573139369488
$ws.Cells.Item(26, 2).Value = 0.9959460496902466
$ws.Cells.Item(26, 3).Value = 0.0001041729483404197
$ws.Cells.Item(26, 4).Value = 1

$ws.Cells.Item(27, 1).Value = 0.009497462771832943
$ws.Cells.Item(27, 2).Value = 0.9962069392204285
$ws.Cells.Item(27, 3).Value = 0.00004726613769889809
$ws.Cells.Item(27, 4).Value = 1

$ws.Cells.Item(28, 1).Value = 0.009835576638579369
$ws.Cells.Item(28, 2).Value = 0.996026337146759
$ws.Cells.Item(28, 3).Value = 0.00003178512997692451
$ws.Cells.Item(28, 4).Value = 1

$ws.Cells.Item(29, 1).Value = 0.009828174486756325
$ws.Cells.Item(29, 2).Value = 0.9962069392204285
$ws.Cells.Item(29, 3).Value = 0.00003441358057898469
$ws.Cells.Item(29, 4).Value = 1

$ws.Cells.Item(30, 1).Value = 0.01051153335720301
$ws.Cells.Item(30, 2).Value = 0.9959861636161804
$ws.Cells.Item(30, 3).Value = 0.00004798976078745909
$ws.Cells.Item(30, 4).Value = 1

$ws.Cells.Item(31, 1).Value = 0.01000656839460135
$ws.Cells.Item(31, 2).Value = 0.9959661364555359
$ws.Cells.Item(31, 3).Value = 0.000056865163060138
$ws.Cells.Item(31, 4).Value = 1

$ws.Cells.Item(32, 1).Value = 0.01003664825111628
$ws.Cells.Item(32, 2).Value = 0.9959059357643127
$ws.Cells.Item(32, 3).Value = 0.00002886227048293222
$ws.Cells.Item(32, 4).Value = 1

$ws.Cells.Item(33, 1).Value = 0.009769693948328495
$ws.Cells.Item(33, 2).Value = 0.9960664510726929
$ws.Cells.Item(33, 3).Value = 0.00002760546522040386
$ws.Cells.Item(33, 4).Value = 1

$ws.Cells.Item(34, 1).Value = 0.009829587303102016
$ws.Cells.Item(34, 2).Value = 0.9959861636161804
$ws.Cells.Item(34, 3).Value = 0.00001194068499899004
$ws.Cells.Item(34, 4).Value = 1

$ws.Cells.Item(35, 1).Value = 0.009488164447247982
$ws.Cells.Item(35, 2).Value = 0.9962471127510071
$ws.Cells.Item(35, 3).Value = 0.0000179078706423752
$ws.Cells.Item(35, 4).Value = 1

$ws.Cells.Item(36, 1).Value = 0.009228608570992947
$ws.Cells.Item(36, 2).Value = 0.9964277148246765
$ws.Cells.Item(36, 3).Value = 0.00003025570913450792
$ws.Cells.Item(36, 4).Value = 1

$ws.Cells.Item(37, 1).Value = 0.009703014045953751
$ws.Cells.Item(37, 2).Value = 0.9961668252944946
$ws.Cells.Item(37, 3).Value = 0.0002504869480617344
$ws.Cells.Item(37, 4).Value = 0.9998840689659119

$ws.Cells.Item(38, 1).Value = 0.009812811389565468
$ws.Cells.Item(38, 2).Value = 0.9962069392204285
$ws.Cells.Item(38, 3).Value = 0.00004119776713196188
$ws.Cells.Item(38, 4).Value = 1

$ws.Cells.Item(39, 1).Value = 0.01000701449811459
$ws.Cells.Item(39, 2).Value = 0.9964478015899658
$ws.Cells.Item(39, 3).Value = 0.00005411657912191004
$ws.Cells.Item(39, 4).Value = 1

$ws.Cells.Item(40, 1).Value = 0.01028370950371027
$ws.Cells.Item(40, 2).Value = 0.9963876008987427
$ws.Cells.Item(40, 3).Value = 0.00004741046723211184
$ws.Cells.Item(40, 4).Value = 1

$ws.Cells.Item(41, 1).Value = 0.00995064340531826
$ws.Cells.Item(41, 2).Value = 0.996126651763916
$ws.Cells.Item(41, 3).Value = 0.0000834599559311755
$ws.Cells.Item(41, 4).Value = 1

$ws.Cells.Item(42, 1).Value = 0.009548894129693508
$ws.Cells.Item(42, 2).Value = 0.9962671399116516
$ws.Cells.Item(42, 3).Value = 0.0000957400188781321
$ws.Cells.Item(42, 4).Value = 1

$ws.Cells.Item(43, 1).Value = 0.01061257533729076
$ws.Cells.Item(43, 2).Value = 0.9956249594688416
$ws.Cells.Item(43, 3).Value = 0.0001257478288607672
$ws.Cells.Item(43, 4).Value = 1

$ws.Cells.Item(44, 1).Value = 0.00897060614079237
$ws.Cells.Item(44, 2).Value = 0.9966284036636353
$ws.Cells.Item(44, 3).Value = 0.0003120446635875851
$ws.Cells.Item(44, 4).Value = 0.9998840689659119

$ws.Cells.Item(45, 1).Value = 0.009908711537718773
$ws.Cells.Item(45, 2).Value = 0.996026337146759
$ws.Cells.Item(45, 3).Value = 0.00036870181793347
$ws.Cells.Item(45, 4).Value = 0.9998840689659119

$ws.Cells.Item(46, 1).Value = 0.009123586118221283
$ws.Cells.Item(46, 2).Value = 0.9964478015899658
$ws.Cells.Item(46, 3).Value = 0.00005266796870273538
$ws.Cells.Item(46, 4).Value = 1

$ws.Cells.Item(47, 1).Value = 0.009714526124298573
$ws.Cells.Item(47, 2).Value = 0.9960865378379822
$ws.Cells.Item(47, 3).Value = 0.00007479038322344422
$ws.Cells.Item(47, 4).Value = 1

$ws.Cells.Item(48, 1).Value = 0.009951287880539894
$ws.Cells.Item(48, 2).Value = 0.9959661364555359
$ws.Cells.Item(48, 3).Value = 0.0005127583863213658
$ws.Cells.Item(48, 4).Value = 0.9998840689659119

$ws.Cells.Item(49, 1).Value = 0.008780322968959808
$ws.Cells.Item(49, 2).Value = 0.9966886043548584
$ws.Cells.Item(49, 3).Value = 0.00001580136631673668
$ws.Cells.Item(49, 4).Value = 1

$ws.Cells.Item(50, 1).Value = 0.01017113029956818
$ws.Cells.Item(50, 2).Value = 0.9959259629249573
$ws.Cells.Item(50, 3).Value = 0.00001238000822922913
$ws.Cells.Item(50, 4).Value = 1

$ws.Cells.Item(51, 1).Value = 0.009123872965574265
$ws.Cells.Item(51, 2).Value = 0.9963675141334534
$ws.Cells.Item(51, 3).Value = 0.00001171946132672019
$ws.Cells.Item(51, 4).Value = 1

